# Applies the 2023-07-03 cryptos list refresh: updates Price/Volume(1h)
# figures for every coin row and fixes the ShibaInu / WrappedliquidstakedEther2.0
# row ordering (rows 20-21) including their Coin name and Link columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.675.93"
$ws.Range("E2").Value = "'  +0.54%  "
$ws.Range("D3").Value = "'1.962.23"
$ws.Range("E3").Value = "'  +2.45%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'248.78"
$ws.Range("E5").Value = "'  +1.38%  "
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("D7").Value = "'0.4838"
$ws.Range("E7").Value = "'  +0.67%  "
$ws.Range("D8").Value = "'0.2954"
$ws.Range("E8").Value = "'  +2.50%  "
$ws.Range("D9").Value = "'0.06790"
$ws.Range("E9").Value = "'  +0.89%  "
$ws.Range("D10").Value = "'110.93"
$ws.Range("E10").Value = "'  +0.09%  "
$ws.Range("E11").Value = "'  +0.94%  "
$ws.Range("D12").Value = "'1.966.04"
$ws.Range("E12").Value = "'  +2.71%  "
$ws.Range("D13").Value = "'0.07741"
$ws.Range("E13").Value = "'  +2.43%  "
$ws.Range("D14").Value = "'5.493"
$ws.Range("E14").Value = "'  +4.62%  "
$ws.Range("D15").Value = "'0.6913"
$ws.Range("E15").Value = "'  +3.35%  "
$ws.Range("D16").Value = "'294.08"
$ws.Range("E16").Value = "'  +0.39%  "
$ws.Range("D17").Value = "'30.694.53"
$ws.Range("E17").Value = "'  +0.63%  "
$ws.Range("D18").Value = "'13.31"
$ws.Range("E18").Value = "'  +2.94%  "
$ws.Range("D19").Value = "'5.674"
$ws.Range("E19").Value = "'  +3.51%  "
$ws.Range("B20").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "'2.226.66"
$ws.Range("E20").Value = "'  +2.93%  "
$ws.Range("B21").Value = "'ShibaInu"
$ws.Range("C21").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000007690"
$ws.Range("E21").Value = "'  +1.44%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "'  +0.03%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "'  +0.06%  "
$ws.Range("D24").Value = "'6.612"
$ws.Range("E24").Value = "'  +3.20%  "
$ws.Range("D25").Value = "'9.934"
$ws.Range("E25").Value = "'  +4.85%  "
$ws.Range("D26").Value = "'170.77"
$ws.Range("E26").Value = "'  +3.94%  "
$ws.Range("D27").Value = "'20.16"
$ws.Range("E27").Value = "'  -1.48%  "
$ws.Range("E28").Value = "'  +3.82%  "
$ws.Range("D29").Value = "'0.1073"
$ws.Range("E29").Value = "'  +0.37%  "
$ws.Range("D30").Value = "'1.443"
$ws.Range("E30").Value = "'  +2.88%  "
$ws.Range("D31").Value = "'4.698"
$ws.Range("E31").Value = "'  +16.68%  "
$ws.Range("D32").Value = "'4.471"
$ws.Range("E32").Value = "'  +7.28%  "
$ws.Range("D33").Value = "'0.05138"
$ws.Range("E33").Value = "'  +3.34%  "
$ws.Range("D34").Value = "'0.7798"
$ws.Range("E34").Value = "'  +6.81%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "'  +3.94%  "
$ws.Range("D36").Value = "'0.02062"
$ws.Range("E36").Value = "'  +0.28%  "
$ws.Range("D37").Value = "'2.735"
$ws.Range("E37").Value = "'  +0.07%  "
$ws.Range("D38").Value = "'2.712"
$ws.Range("D39").Value = "'2.072"
$ws.Range("E39").Value = "'  +2.86%  "
$ws.Range("D40").Value = "'111.23"
$ws.Range("E40").Value = "'  +0.25%  "
$ws.Range("D41").Value = "'6.136"
$ws.Range("E41").Value = "'  +4.45%  "
$ws.Range("D42").Value = "'0.4468"
$ws.Range("E42").Value = "'  +1.04%  "
$ws.Range("D43").Value = "'0.8748"
$ws.Range("E43").Value = "'  +1.21%  "
$ws.Range("D44").Value = "'70.32"
$ws.Range("E44").Value = "'  +3.05%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("D46").Value = "'7.393"
$ws.Range("E46").Value = "'  +1.05%  "
$ws.Range("D47").Value = "'0.1280"
$ws.Range("E47").Value = "'  +3.58%  "
$ws.Range("D48").Value = "'9.418"
$ws.Range("E48").Value = "'  +1.28%  "
$ws.Range("E49").Value = "'  +3.18%  "
$ws.Range("D50").Value = "'47.77"
$ws.Range("E50").Value = "'  -2.34%  "
$ws.Range("D51").Value = "'0.2516"
$ws.Range("E51").Value = "'  -0.92%  "
